$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the last four table columns (header cells in row 2) to be prefixed with "1_"
$ws.Range("I2").Value = "1_50"
$ws.Range("J2").Value = "1_100"
$ws.Range("K2").Value = "1_150"
$ws.Range("L2").Value = "1_latest"

# The wider header text makes Excel recompute the (best-fit) column widths
$ws.Columns.Item(9).ColumnWidth = 6.333333333333333
$ws.Range("J1:K1").EntireColumn.ColumnWidth = 7.333333333333333
$ws.Columns.Item(12).ColumnWidth = 9

# Move/record the active selection as it was when the workbook was saved
$ws.Range("K4").Select() | Out-Null
